# 2018-02-07 First Event Updates
# Edits slide 37 ("A developer creates an application ...") :
#   - Title: bold+blue ("0070C0") highlight on "an application" and "external resources"
#   - Content placeholder: bold+blue highlight on "create", plus an explicit
#     position/size (xfrm) on the placeholder shape
#   - Two "Rectangle" callout boxes get repositioned / resized

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(37)

# ---------------------------------------------------------------------------
# Shape 1 : Title 7
# ---------------------------------------------------------------------------
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange

# "A developer creates " (1-20) / "an application " (21-35, bold+blue) /
# "that needs access to " (36-56) / "external resources" (57-74, bold+blue) /
# ".  This application will be deployed in the domain.  " (75-127)
$run1 = $titleRange.Characters(21, 15)
$run1.Font.Bold = $true
$run1.Font.Color.RGB = 12611584

$run2 = $titleRange.Characters(57, 18)
$run2.Font.Bold = $true
$run2.Font.Color.RGB = 12611584

# ---------------------------------------------------------------------------
# Shape 2 : Content Placeholder 8
# ---------------------------------------------------------------------------
$content = $s.Shapes.Item(2)

# explicit size/position (previously inherited from the layout)
$content.Left = 20.565906524658203
$content.Top = 162.7076416015625
$content.Width = 675.1840209960938
$content.Height = 325.51055908203125

$contentRange = $content.TextFrame.TextRange

# "You need to use Azure Command-Line Interface (CLI) to " (1-54) /
# "create " (55-61, bold+blue) /
# "a service principle.  How should you configure the command? " (62-121)
$run3 = $contentRange.Characters(55, 7)
$run3.Font.Bold = $true
$run3.Font.Color.RGB = 12611584

# ---------------------------------------------------------------------------
# Shape 4 : Rectangle 10
# ---------------------------------------------------------------------------
$rect10 = $s.Shapes.Item(4)
$rect10.Left = 376.84307861328125
$rect10.Top = 324.0948181152344
$rect10.Width = 152.69378662109375
$rect10.Height = 25.116300582885742

# ---------------------------------------------------------------------------
# Shape 5 : Rectangle 11
# ---------------------------------------------------------------------------
$rect11 = $s.Shapes.Item(5)
$rect11.Left = 375.0697937011719
$rect11.Top = 425.2226867675781
$rect11.Width = 135.6278839111328
$rect11.Height = 25.116300582885742
